$d = $word.ActiveDocument

# Locate the "Pick up prompt ... walls" bullet; everything else anchors
# off this paragraph's index so the script isn't brittle to unrelated
# paragraph-count assumptions.
$pickUpIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*Pick up prompt*walls*") {
        $pickUpIndex = $i
        break
    }
}

# ------------------------------------------------------------------
# 1) Wrap the last word ("walls") of the "Pick up prompt..." bullet in
#    a grammar-checker proofErr pair, splitting the run in two.
# ------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("walls", $false, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
if ($found) {
    # Re-seat the hit as a plain Range (the live Find range behaves like an
    # insertion point for InsertXML) so InsertXML replaces the word itself
    # rather than appending after it.
    $wallsRng = $d.Range($findRng.Start, $findRng.End)
    $wallsXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' +
        '<w:p>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>walls</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $wallsRng.InsertXML($wallsXml)
}

# ------------------------------------------------------------------
# 2) Replace the two trailing empty paragraphs (immediately after the
#    "Pick up prompt..." bullet) with a new bulleted list item:
#    "Texture Streaming Size to 2000Mb" (three runs).
# ------------------------------------------------------------------
$pLast = $d.Paragraphs.Item($pickUpIndex + 2)
$pPrev = $d.Paragraphs.Item($pickUpIndex + 1)

# Merge the two trailing empty paragraphs into a single empty paragraph.
$mergeRng = $d.Range($pPrev.Range.Start, $pLast.Range.End)
$mergeRng.Delete()

# Fill that single empty paragraph with the fully-formatted new list
# item (style + numbering + three text runs) via InsertXML.
$target = $d.Paragraphs.Item($d.Paragraphs.Count)
$newParaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:pPr>' +
    '<w:pStyle w:val="ListParagraph"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr>' +
    '</w:pPr>' +
    '<w:r><w:t xml:space="preserve">Texture Streaming Size to </w:t></w:r>' +
    '<w:r><w:t>2000</w:t></w:r>' +
    '<w:r><w:t>Mb</w:t></w:r>' +
    '</w:p>' +
    '</w:body></w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$target.Range.InsertXML($newParaXml)

# InsertXML above inserted the new paragraph in front of the old empty
# paragraph mark, leaving that now-redundant empty paragraph trailing
# at the very end of the document body; remove it by merging it back
# into the paragraph we just created.
$newTextPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$trailingEmpty = $d.Paragraphs.Item($d.Paragraphs.Count)
$cleanupRng = $d.Range($newTextPara.Range.End - 1, $trailingEmpty.Range.End)
$cleanupRng.Delete()

# ------------------------------------------------------------------
# 3) Update the "Last Updated" timestamp in the footer.
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)
$footer.Range.Find.Execute("19/01/2023 16:45", $false, $false, $false, $false, `
                            $false, $true, 1, $false, "19/01/2023 23:28", 2)
